$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2,1).Value = 'Última actualización: 11:55:01'
$ws.Cells.Item(3,1).Value = 'Total filas: 165'
$ws.Cells.Item(16,1).Value = '06:38:54'
$ws.Cells.Item(16,3).Value = '16_SANTA ANA'
$ws.Cells.Item(16,4).Value = 2
$ws.Cells.Item(17,1).Value = '05:44:02'
$ws.Cells.Item(17,3).Value = '17X38_ROMERO'
$ws.Cells.Item(17,4).Value = 56
$ws.Cells.Item(28,1).Value = '06:56:24'
$ws.Cells.Item(28,3).Value = '16_SANTA ANA'
$ws.Cells.Item(28,4).Value = 25
$ws.Cells.Item(29,1).Value = '07:15:48'
$ws.Cells.Item(29,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(29,4).Value = 6
$ws.Cells.Item(49,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(50,3).Value = '15_ABASTO'
$ws.Cells.Item(51,3).Value = '15_ABASTO'
$ws.Cells.Item(52,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(55,1).Value = '08:30:14'
$ws.Cells.Item(55,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(55,4).Value = 11
$ws.Cells.Item(56,1).Value = '07:52:32'
$ws.Cells.Item(56,3).Value = '10_OLMOS'
$ws.Cells.Item(56,4).Value = 49
$ws.Cells.Item(64,3).Value = '215B_EL PATO'
$ws.Cells.Item(65,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(75,1).Value = '08:30:14'
$ws.Cells.Item(75,3).Value = '27_EL RETIRO'
$ws.Cells.Item(75,4).Value = 47
$ws.Cells.Item(76,1).Value = '08:40:59'
$ws.Cells.Item(76,3).Value = '15X38_ABASTO'
$ws.Cells.Item(76,4).Value = 37
$ws.Cells.Item(77,1).Value = '08:52:33'
$ws.Cells.Item(77,3).Value = '14_ABASTO'
$ws.Cells.Item(77,4).Value = 25
$ws.Cells.Item(78,1).Value = '08:52:33'
$ws.Cells.Item(78,3).Value = '15X38_ABASTO'
$ws.Cells.Item(78,4).Value = 26
$ws.Cells.Item(79,1).Value = '08:30:14'
$ws.Cells.Item(79,3).Value = '14_ABASTO'
$ws.Cells.Item(79,4).Value = 48
$ws.Cells.Item(112,1).Value = '10:56:01'
$ws.Cells.Item(112,3).Value = '81_EL PELIGRO'
$ws.Cells.Item(112,4).Value = 5
$ws.Cells.Item(113,1).Value = '09:23:52'
$ws.Cells.Item(113,3).Value = '10_OLMOS'
$ws.Cells.Item(113,4).Value = 98
$ws.Cells.Item(134,1).Value = '11:48:20'
$ws.Cells.Item(134,3).Value = '225_GOMEZ'
$ws.Cells.Item(134,4).Value = 5
$ws.Cells.Item(135,1).Value = '11:35:40'
$ws.Cells.Item(135,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(135,4).Value = 18
$ws.Cells.Item(136,1).Value = '11:55:01'
$ws.Cells.Item(136,2).Value = '11:55'
$ws.Cells.Item(136,3).Value = '225_GOMEZ'
$ws.Cells.Item(136,4).Value = 0
$ws.Cells.Item(137,1).Value = '11:55:01'
$ws.Cells.Item(137,2).Value = '11:58'
$ws.Cells.Item(137,3).Value = '17_ROMERO'
$ws.Cells.Item(137,4).Value = 3
$ws.Cells.Item(138,1).Value = '11:35:40'
$ws.Cells.Item(138,2).Value = '12:05'
$ws.Cells.Item(138,4).Value = 30
$ws.Cells.Item(139,1).Value = '11:13:01'
$ws.Cells.Item(139,2).Value = '12:06'
$ws.Cells.Item(139,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(139,4).Value = 53
$ws.Cells.Item(140,1).Value = '11:55:01'
$ws.Cells.Item(140,4).Value = 15
$ws.Cells.Item(141,1).Value = '11:55:01'
$ws.Cells.Item(141,2).Value = '12:10'
$ws.Cells.Item(141,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(141,4).Value = 15
$ws.Cells.Item(142,1).Value = '11:55:01'
$ws.Cells.Item(142,2).Value = '12:17'
$ws.Cells.Item(142,3).Value = '10_OLMOS'
$ws.Cells.Item(142,4).Value = 22
$ws.Cells.Item(143,1).Value = '11:35:40'
$ws.Cells.Item(143,2).Value = '12:21'
$ws.Cells.Item(143,4).Value = 46
$ws.Cells.Item(144,1).Value = '11:55:01'
$ws.Cells.Item(144,2).Value = '12:22'
$ws.Cells.Item(144,3).Value = '215C_EL PATO'
$ws.Cells.Item(144,4).Value = 27
$ws.Cells.Item(145,1).Value = '11:13:01'
$ws.Cells.Item(145,3).Value = '27_EL RETIRO'
$ws.Cells.Item(145,4).Value = 78
$ws.Cells.Item(146,1).Value = '11:35:40'
$ws.Cells.Item(146,2).Value = '12:31'
$ws.Cells.Item(146,4).Value = 56
$ws.Cells.Item(148,1).Value = '11:48:20'
$ws.Cells.Item(148,2).Value = '12:32'
$ws.Cells.Item(148,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(148,4).Value = 44
$ws.Cells.Item(150,1).Value = '11:55:01'
$ws.Cells.Item(150,2).Value = '12:33'
$ws.Cells.Item(150,3).Value = '14_ABASTO'
$ws.Cells.Item(150,4).Value = 38
$ws.Cells.Item(152,1).Value = '11:55:01'
$ws.Cells.Item(152,2).Value = '12:34'
$ws.Cells.Item(152,3).Value = '15_ABASTO'
$ws.Cells.Item(152,4).Value = 39
$ws.Cells.Item(153,1).Value = '11:35:40'
$ws.Cells.Item(153,2).Value = '12:36'
$ws.Cells.Item(153,4).Value = 61
$ws.Cells.Item(154,1).Value = '11:55:01'
$ws.Cells.Item(154,2).Value = '12:37'
$ws.Cells.Item(154,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(154,4).Value = 42
$ws.Cells.Item(155,1).Value = '11:55:01'
$ws.Cells.Item(155,2).Value = '12:37'
$ws.Cells.Item(155,3).Value = '27_EL RETIRO'
$ws.Cells.Item(155,4).Value = 42
$ws.Cells.Item(156,2).Value = '12:47'
$ws.Cells.Item(156,3).Value = '14_ABASTO'
$ws.Cells.Item(156,4).Value = 59
$ws.Cells.Item(157,1).Value = '11:55:01'
$ws.Cells.Item(157,2).Value = '12:48'
$ws.Cells.Item(157,3).Value = '14_ABASTO'
$ws.Cells.Item(157,4).Value = 53
$ws.Cells.Item(158,1).Value = '11:55:01'
$ws.Cells.Item(158,2).Value = '12:48'
$ws.Cells.Item(158,3).Value = '16_SANTA ANA'
$ws.Cells.Item(158,4).Value = 53
$ws.Cells.Item(159,1).Value = '11:55:01'
$ws.Cells.Item(159,2).Value = '12:48'
$ws.Cells.Item(159,3).Value = '15X38_ABASTO'
$ws.Cells.Item(159,4).Value = 53
$ws.Cells.Item(160,2).Value = '13:02'
$ws.Cells.Item(160,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(160,4).Value = 74
$ws.Cells.Item(161,1).Value = '11:55:01'
$ws.Cells.Item(161,2).Value = '13:03'
$ws.Cells.Item(161,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(161,4).Value = 68
$ws.Cells.Item(162,1).Value = '11:35:40'
$ws.Cells.Item(162,2).Value = '13:03'
$ws.Cells.Item(162,3).Value = '215C_EL PATO'
$ws.Cells.Item(162,4).Value = 88
$ws.Cells.Item(163,1).Value = '11:55:01'
$ws.Cells.Item(163,2).Value = '13:04'
$ws.Cells.Item(163,3).Value = '215C_EL PATO'
$ws.Cells.Item(163,4).Value = 69
$ws.Cells.Item(164,1).Value = '11:55:01'
$ws.Cells.Item(164,2).Value = '13:13'
$ws.Cells.Item(164,3).Value = '16_SANTA ANA'
$ws.Cells.Item(164,4).Value = 78
$ws.Cells.Item(165,1).Value = '11:55:01'
$ws.Cells.Item(165,2).Value = '13:17'
$ws.Cells.Item(165,3).Value = '10_OLMOS'
$ws.Cells.Item(165,4).Value = 82
$ws.Cells.Item(165,5).Value = 'LP1912'
$ws.Cells.Item(166,1).Value = '11:55:01'
$ws.Cells.Item(166,2).Value = '13:24'
$ws.Cells.Item(166,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(166,4).Value = 89
$ws.Cells.Item(166,5).Value = 'LP1912'
$ws.Cells.Item(167,1).Value = '11:55:01'
$ws.Cells.Item(167,2).Value = '13:25'
$ws.Cells.Item(167,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(167,4).Value = 90
$ws.Cells.Item(167,5).Value = 'LP1912'
$ws.Cells.Item(168,1).Value = '11:55:01'
$ws.Cells.Item(168,2).Value = '13:33'
$ws.Cells.Item(168,3).Value = '215A_EL PATO'
$ws.Cells.Item(168,4).Value = 98
$ws.Cells.Item(168,5).Value = 'LP1912'
$ws.Cells.Item(169,1).Value = '11:55:01'
$ws.Cells.Item(169,2).Value = '13:47'
$ws.Cells.Item(169,3).Value = '225_GOMEZ'
$ws.Cells.Item(169,4).Value = 112
$ws.Cells.Item(169,5).Value = 'LP1912'
$ws.Cells.Item(170,1).Value = '11:55:01'
$ws.Cells.Item(170,2).Value = '13:49'
$ws.Cells.Item(170,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(170,4).Value = 114
$ws.Cells.Item(170,5).Value = 'LP1912'

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2,1).Value = 'Última actualización: 11:55:01'
$ws.Cells.Item(26,1).Value = '11:55:01'
$ws.Cells.Item(26,4).Value = 27
$ws.Cells.Item(28,1).Value = '11:55:01'
$ws.Cells.Item(28,4).Value = 69
$ws.Cells.Item(29,1).Value = '11:55:01'
$ws.Cells.Item(29,4).Value = 98

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2,1).Value = 'Última actualización: 11:55:01'
$ws.Cells.Item(25,1).Value = '11:55:01'
$ws.Cells.Item(25,4).Value = 77
$ws.Cells.Item(27,1).Value = '11:55:01'
$ws.Cells.Item(27,4).Value = 86
